# Slide 12 ("M16: Het project gebruikt tools voor vastgestelde taken"):
# combine the two separate bullet lists (tasks + tools) that used to live
# in "TextBox 2" into a single explanatory paragraph.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item("TextBox 2")
$shape.TextFrame.TextRange.Text = "Voor vastgestelde taken bij het ontwikkelen, onderhouden en operationeel beheren van software, stelt ICTU het gebruik van tools verplicht. ICTU adviseert per taak specifieke tools en ondersteunt projecten bij het gebruik daarvan."
